$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values for re-annotated rows
$ws.Cells.Item(2, 9).Value = "sv"
$ws.Cells.Item(2, 10).Value = "Statement-opinion"
$ws.Cells.Item(7, 9).Value = "sv"
$ws.Cells.Item(7, 10).Value = "Statement-opinion"
$ws.Cells.Item(11, 9).Value = "sd"
$ws.Cells.Item(11, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(18, 9).Value = "aa"
$ws.Cells.Item(18, 10).Value = "Agree/Accept"
$ws.Cells.Item(42, 9).Value = "sd"
$ws.Cells.Item(42, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(50, 9).Value = "sd"
$ws.Cells.Item(50, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(57, 9).Value = "%"
$ws.Cells.Item(57, 10).Value = "Uninterpretable"
$ws.Cells.Item(89, 9).Value = "sd"
$ws.Cells.Item(89, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(98, 9).Value = "b"
$ws.Cells.Item(98, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(102, 9).Value = "sd"
$ws.Cells.Item(102, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(110, 9).Value = "aa"
$ws.Cells.Item(110, 10).Value = "Agree/Accept"
$ws.Cells.Item(114, 9).Value = "sv"
$ws.Cells.Item(114, 10).Value = "Statement-opinion"
$ws.Cells.Item(118, 9).Value = "sd"
$ws.Cells.Item(118, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(123, 9).Value = "sv"
$ws.Cells.Item(123, 10).Value = "Statement-opinion"
$ws.Cells.Item(128, 9).Value = "b"
$ws.Cells.Item(128, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(144, 9).Value = "aa"
$ws.Cells.Item(144, 10).Value = "Agree/Accept"
$ws.Cells.Item(147, 9).Value = "sd"
$ws.Cells.Item(147, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(151, 9).Value = "sd"
$ws.Cells.Item(151, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(157, 9).Value = "sd"
$ws.Cells.Item(157, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(158, 9).Value = "sd"
$ws.Cells.Item(158, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(174, 9).Value = "aa"
$ws.Cells.Item(174, 10).Value = "Agree/Accept"
$ws.Cells.Item(178, 9).Value = "sd"
$ws.Cells.Item(178, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(181, 9).Value = "sv"
$ws.Cells.Item(181, 10).Value = "Statement-opinion"
$ws.Cells.Item(189, 9).Value = "sd"
$ws.Cells.Item(189, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(190, 9).Value = "sd"
$ws.Cells.Item(190, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(191, 9).Value = "sd"
$ws.Cells.Item(191, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(194, 9).Value = "sv"
$ws.Cells.Item(194, 10).Value = "Statement-opinion"
$ws.Cells.Item(207, 9).Value = "sd"
$ws.Cells.Item(207, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(208, 9).Value = "sv"
$ws.Cells.Item(208, 10).Value = "Statement-opinion"
$ws.Cells.Item(209, 9).Value = "sd"
$ws.Cells.Item(209, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(213, 9).Value = "sv"
$ws.Cells.Item(213, 10).Value = "Statement-opinion"
$ws.Cells.Item(214, 9).Value = "sd"
$ws.Cells.Item(214, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(215, 9).Value = "sd"
$ws.Cells.Item(215, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(220, 9).Value = "sd"
$ws.Cells.Item(220, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(221, 9).Value = "sv"
$ws.Cells.Item(221, 10).Value = "Statement-opinion"
$ws.Cells.Item(222, 9).Value = "sv"
$ws.Cells.Item(222, 10).Value = "Statement-opinion"
$ws.Cells.Item(224, 9).Value = "sv"
$ws.Cells.Item(224, 10).Value = "Statement-opinion"
$ws.Cells.Item(230, 9).Value = "qy"
$ws.Cells.Item(230, 10).Value = "Yes-No-Question"
$ws.Cells.Item(232, 9).Value = "b"
$ws.Cells.Item(232, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(233, 9).Value = "sd"
$ws.Cells.Item(233, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(261, 9).Value = "sv"
$ws.Cells.Item(261, 10).Value = "Statement-opinion"
$ws.Cells.Item(263, 9).Value = "sv"
$ws.Cells.Item(263, 10).Value = "Statement-opinion"
$ws.Cells.Item(277, 9).Value = "sv"
$ws.Cells.Item(277, 10).Value = "Statement-opinion"
$ws.Cells.Item(298, 9).Value = "aa"
$ws.Cells.Item(298, 10).Value = "Agree/Accept"
$ws.Cells.Item(299, 9).Value = "aa"
$ws.Cells.Item(299, 10).Value = "Agree/Accept"
$ws.Cells.Item(312, 9).Value = "sd"
$ws.Cells.Item(312, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(314, 9).Value = "%"
$ws.Cells.Item(314, 10).Value = "Uninterpretable"
$ws.Cells.Item(317, 9).Value = "sd"
$ws.Cells.Item(317, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(322, 9).Value = "%"
$ws.Cells.Item(322, 10).Value = "Uninterpretable"
$ws.Cells.Item(325, 9).Value = "aa"
$ws.Cells.Item(325, 10).Value = "Agree/Accept"
$ws.Cells.Item(332, 9).Value = "sd"
$ws.Cells.Item(332, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(335, 9).Value = "sd"
$ws.Cells.Item(335, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(340, 9).Value = "sv"
$ws.Cells.Item(340, 10).Value = "Statement-opinion"
$ws.Cells.Item(343, 9).Value = "ba"
$ws.Cells.Item(343, 10).Value = "Appreciation"
$ws.Cells.Item(361, 9).Value = "aa"
$ws.Cells.Item(361, 10).Value = "Agree/Accept"
$ws.Cells.Item(363, 9).Value = "sd"
$ws.Cells.Item(363, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(368, 9).Value = "aa"
$ws.Cells.Item(368, 10).Value = "Agree/Accept"
$ws.Cells.Item(369, 9).Value = "aa"
$ws.Cells.Item(369, 10).Value = "Agree/Accept"
$ws.Cells.Item(374, 9).Value = "ba"
$ws.Cells.Item(374, 10).Value = "Appreciation"
$ws.Cells.Item(390, 9).Value = "%"
$ws.Cells.Item(390, 10).Value = "Uninterpretable"
$ws.Cells.Item(400, 9).Value = "sd"
$ws.Cells.Item(400, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(402, 9).Value = "aa"
$ws.Cells.Item(402, 10).Value = "Agree/Accept"
$ws.Cells.Item(406, 9).Value = "ba"
$ws.Cells.Item(406, 10).Value = "Appreciation"
